$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 980.5465
$ws.Range("I129").Value = 743.5
$ws.Range("J129").Value = 986.1905
$ws.Range("K129").Value = 2230.5
$ws.Range("L129").Value = 2958.5715
$ws.Range("M129").Value = 2769.5
$ws.Range("N129").Value = -12958.5715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2503.0278
$ws.Range("I61").Value = 1479.5
$ws.Range("J61").Value = 5164.2
$ws.Range("K61").Value = 1479.5
$ws.Range("L61").Value = 5164.2
$ws.Range("M61").Value = -1267.5
$ws.Range("N61").Value = -5588.2
$ws.Range("H74").Value = 7047.091
$ws.Range("I74").Value = 1774.2222
$ws.Range("J74").Value = 30775
$ws.Range("K74").Value = 1774.2222
$ws.Range("L74").Value = 30775
$ws.Range("M74").Value = -900.2221999999999
$ws.Range("N74").Value = -32523
$ws.Range("H77").Value = 7047.091
$ws.Range("I77").Value = 1774.2222
$ws.Range("J77").Value = 30775
$ws.Range("K77").Value = 8871.110999999999
$ws.Range("L77").Value = 153875
$ws.Range("M77").Value = -4503.110999999999
$ws.Range("N77").Value = -162611
$ws.Range("H132").Value = 2307.0881
$ws.Range("I132").Value = 1771.4642
$ws.Range("K132").Value = 5314.392599999999
$ws.Range("M132").Value = -2784.392599999999
$ws.Range("H136").Value = 2503.0278
$ws.Range("I136").Value = 1479.5
$ws.Range("J136").Value = 5164.2
$ws.Range("K136").Value = 4438.5
$ws.Range("L136").Value = 15492.6
$ws.Range("M136").Value = -1888.5
$ws.Range("N136").Value = -20592.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5812.28
$ws.Range("I86").Value = 1625.75
$ws.Range("J86").Value = 9676.77
$ws.Range("K86").Value = 1625.75
$ws.Range("L86").Value = 9676.77
$ws.Range("M86").Value = -502.75
$ws.Range("N86").Value = -11922.77
$ws.Range("H89").Value = 5812.28
$ws.Range("I89").Value = 1625.75
$ws.Range("J89").Value = 9676.77
$ws.Range("K89").Value = 8128.75
$ws.Range("L89").Value = 48383.85000000001
$ws.Range("M89").Value = -2512.75
$ws.Range("N89").Value = -59615.85000000001
$ws.Range("H105").Value = 7582.2383
$ws.Range("I105").Value = 8431
$ws.Range("J105").Value = 3975
$ws.Range("K105").Value = 8431
$ws.Range("L105").Value = 3975
$ws.Range("M105").Value = -6684
$ws.Range("N105").Value = -7469
$ws.Range("H134").Value = 3919.375
$ws.Range("I134").Value = 3021.5334
$ws.Range("J134").Value = 5415.778
$ws.Range("K134").Value = 9064.600199999999
$ws.Range("L134").Value = 16247.334
$ws.Range("M134").Value = -6529.600199999999
$ws.Range("N134").Value = -21317.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 125327.75
$ws.Range("I16").Value = 125327.75
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 125327.75
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -125040.75
$ws.Range("H88").Value = 19638.166
$ws.Range("J88").Value = 19638.166
$ws.Range("L88").Value = 19638.166
$ws.Range("N88").Value = -20450.166
$ws.Range("H91").Value = 19638.166
$ws.Range("J91").Value = 19638.166
$ws.Range("L91").Value = 19638.166
$ws.Range("N91").Value = -22446.166
$ws.Range("H107").Value = 724.25
$ws.Range("I107").Value = 95
$ws.Range("J107").Value = 934
$ws.Range("K107").Value = 95
$ws.Range("L107").Value = 934
$ws.Range("M107").Value = 1825
$ws.Range("N107").Value = -4774
$ws.Range("H113").Value = 125327.75
$ws.Range("I113").Value = 125327.75
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 125327.75
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -123157.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2331
$ws.Range("J75").Value = 980
$ws.Range("L75").Value = 2940
$ws.Range("N75").Value = -4936
$ws.Range("H78").Value = 2331
$ws.Range("J78").Value = 980
$ws.Range("L78").Value = 8820
$ws.Range("N78").Value = -18804
$ws.Range("H113").Value = 15152324
$ws.Range("J113").Value = 16129858
$ws.Range("L113").Value = 48389574
$ws.Range("N113").Value = -48393914
$ws.Range("H131").Value = 1753.4546
$ws.Range("I131").Value = 555.6667
$ws.Range("J131").Value = 1942.579
$ws.Range("K131").Value = 1667.0001
$ws.Range("L131").Value = 5827.737
$ws.Range("M131").Value = 3372.9999
$ws.Range("N131").Value = -15907.737

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1390390.2
$ws.Range("I122").Value = 11111111
$ws.Range("K122").Value = 33333333
$ws.Range("M122").Value = -33330883

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 905.13043
$ws.Range("I16").Value = 920.8
$ws.Range("K16").Value = 920.8
$ws.Range("M16").Value = -750.8
$ws.Range("H68").Value = 1367.5
$ws.Range("I68").Value = 1300.3334
$ws.Range("J68").Value = 1501.8334
$ws.Range("K68").Value = 1300.3334
$ws.Range("L68").Value = 1501.8334
$ws.Range("M68").Value = -551.3334
$ws.Range("N68").Value = -2999.8334
$ws.Range("H71").Value = 1367.5
$ws.Range("I71").Value = 1300.3334
$ws.Range("J71").Value = 1501.8334
$ws.Range("K71").Value = 6501.666999999999
$ws.Range("L71").Value = 7509.166999999999
$ws.Range("M71").Value = -2757.666999999999
$ws.Range("N71").Value = -14997.167
$ws.Range("H93").Value = 543.4375
$ws.Range("I93").Value = 520.11536
$ws.Range("J93").Value = 644.5
$ws.Range("K93").Value = 520.11536
$ws.Range("L93").Value = 644.5
$ws.Range("M93").Value = 727.88464
$ws.Range("N93").Value = -3140.5
$ws.Range("H122").Value = 3598.8948
$ws.Range("I122").Value = 2476
$ws.Range("J122").Value = 3999.9285
$ws.Range("K122").Value = 7428
$ws.Range("L122").Value = 11999.7855
$ws.Range("M122").Value = -4978
$ws.Range("N122").Value = -16899.7855
$ws.Range("H132").Value = 5114.1143
$ws.Range("I132").Value = 2393.818
$ws.Range("J132").Value = 9717.691999999999
$ws.Range("K132").Value = 7181.454000000001
$ws.Range("L132").Value = 29153.076
$ws.Range("M132").Value = -4651.454000000001
$ws.Range("N132").Value = -34213.076
$ws.Range("H136").Value = 6567.2
$ws.Range("I136").Value = 3660.2
$ws.Range("J136").Value = 9474.200000000001
$ws.Range("K136").Value = 10980.6
$ws.Range("L136").Value = 28422.6
$ws.Range("M136").Value = -8430.599999999999
$ws.Range("N136").Value = -33522.60000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 49626
$ws.Range("J108").Value = 49626
$ws.Range("L108").Value = 49626
$ws.Range("N108").Value = -57306
$ws.Range("H122").Value = 54937.105
$ws.Range("I122").Value = 144300.58
$ws.Range("J122").Value = 2808.4167
$ws.Range("K122").Value = 432901.74
$ws.Range("L122").Value = 8425.250100000001
$ws.Range("M122").Value = -430451.74
$ws.Range("N122").Value = -13325.2501
$ws.Range("H132").Value = 14288801
$ws.Range("I132").Value = 27780278
$ws.Range("J132").Value = 3708.5293
$ws.Range("K132").Value = 83340834
$ws.Range("L132").Value = 11125.5879
$ws.Range("M132").Value = -83338304
$ws.Range("N132").Value = -16185.5879
